$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (styles) from column Q (2020) into new columns R, S, T
# for rows 2-5, matching the existing look (borders/number formats/fonts).
$ws.Range("Q2:Q5").Copy()
$ws.Range("R2:T5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new year headers (2021, 2022, 2023)
$ws.Range("R3").Value2 = 2021
$ws.Range("S3").Value2 = 2022
$ws.Range("T3").Value2 = 2023

# Fill in the new data values for "Number of written appeals" row
$ws.Range("R4").Value2 = 4301
$ws.Range("S4").Value2 = 3690
$ws.Range("T4").Value2 = 2620

# Fill in the new data values for "Number of positively resolved" row
$ws.Range("R5").Value2 = 427
$ws.Range("S5").Value2 = 280
$ws.Range("T5").Value2 = 264

# Clear the stale selection in the saved sheet view (it pointed at F16,
# which is outside the used range) so only the tab-selected flag remains.
$ws.Range("A1").Select()
